$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 35
$ws.Range("F3").Value = 27
$ws.Range("H3").Value = 27

$ws.Range("E6").Value = 42

$ws.Range("E10").Value = 448

$ws.Range("E11").Value = 312
$ws.Range("F11").Value = 172
$ws.Range("H11").Value = 172

$ws.Range("E12").Value = 441
$ws.Range("F12").Value = 234
$ws.Range("H12").Value = 234

$ws.Range("E14").Value = 114

$ws.Range("E15").Value = 145
$ws.Range("F15").Value = 58
$ws.Range("H15").Value = 58

$ws.Range("E16").Value = 185
$ws.Range("F16").Value = 92
$ws.Range("H16").Value = 92

$ws.Range("F18").Value = 25
$ws.Range("H18").Value = 25

$ws.Range("E20").Value = 81

$ws.Range("E22").Value = 157
$ws.Range("F22").Value = 78
$ws.Range("H22").Value = 78

$ws.Range("E23").Value = 183

$ws.Range("E24").Value = 193
$ws.Range("F24").Value = 101
$ws.Range("H24").Value = 101

$ws.Range("E25").Value = 238
$ws.Range("F25").Value = 111
$ws.Range("H25").Value = 111

$ws.Range("E26").Value = 143
$ws.Range("F26").Value = 82
$ws.Range("H26").Value = 82

$ws.Range("E27").Value = 300
$ws.Range("F27").Value = 140
$ws.Range("H27").Value = 140

$ws.Range("E28").Value = 183
$ws.Range("F28").Value = 64
$ws.Range("H28").Value = 64

$ws.Range("E29").Value = 155
$ws.Range("F29").Value = 85
$ws.Range("H29").Value = 85

$ws.Range("E30").Value = 192
$ws.Range("F30").Value = 110
$ws.Range("H30").Value = 110

$ws.Range("E32").Value = 171

$ws.Range("F33").Value = 132
$ws.Range("H33").Value = 132

$ws.Range("E34").Value = 199
$ws.Range("F34").Value = 123
$ws.Range("H34").Value = 123

$ws.Range("E36").Value = 64

$ws.Range("E37").Value = 142

$ws.Range("E38").Value = 84
$ws.Range("F38").Value = 54
$ws.Range("H38").Value = 54

$ws.Range("E40").Value = 239

$ws.Range("E41").Value = 363

$ws.Range("E42").Value = 329

$ws.Range("E43").Value = 108

$ws.Range("E44").Value = 291
$ws.Range("F44").Value = 139
$ws.Range("H44").Value = 139

$ws.Range("E45").Value = 128

$ws.Range("E47").Value = 405
$ws.Range("F47").Value = 198
$ws.Range("H47").Value = 198

$ws.Range("E48").Value = 187

$ws.Range("E49").Value = 265
$ws.Range("F49").Value = 112
$ws.Range("H49").Value = 112

$ws.Range("E50").Value = 228

$ws.Range("E51").Value = 217
$ws.Range("F51").Value = 88
$ws.Range("H51").Value = 88

$ws.Range("E52").Value = 24
